$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 05:07:15"
$wsZhCn.Range("H2").Value = "2016-03-22 05:07:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 05:07:18"
$wsDeDe.Range("H2").Value = "2016-03-22 05:08:06"
